# The "江西·东方LiveParty×THO03幻想Strawberry~！！" event (2024-07-13) was
# pulled from the source feed, so it disappears from the "演出" and
# "全部类型" sheets entirely (the "展览" sheet never listed it). Removing
# the row shifts every later row up by one. At the same time the scraper's
# "想去人数" (interest-count) column ticked up for a long list of still-live
# events across all the affected sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. 演出 (Performances): the only data row was the cancelled LiveParty
#    event - delete it, leaving just the header row (A1:I1).
# ---------------------------------------------------------------------------
$wsPerf = $wb.Worksheets.Item("演出")
$wsPerf.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# 2. 全部类型 (All types): same event lived at row 4 here - delete it too,
#    then repair the running index in column A for every row that shifted
#    up (A holds row-number-minus-one).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows.Item(4).Delete()

$lastRow = $wsAll.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $wsAll.Range("A" + $r).Value = $r - 1
}

# ---------------------------------------------------------------------------
# 3. Bump "想去人数" (column F) for both 展览 (already missing the cancelled
#    row) and 全部类型 (now also missing it, after step 2) at the same row
#    numbers.
# ---------------------------------------------------------------------------
$fChanges = @(
    @{Row=2;  New=4911},
    @{Row=3;  New=135},
    @{Row=4;  New=109},
    @{Row=5;  New=801},
    @{Row=6;  New=247},
    @{Row=7;  New=1279},
    @{Row=8;  New=136},
    @{Row=11; New=92},
    @{Row=12; New=7},
    @{Row=13; New=166},
    @{Row=15; New=4309},
    @{Row=16; New=6607},
    @{Row=20; New=555},
    @{Row=22; New=4067},
    @{Row=23; New=427},
    @{Row=24; New=59},
    @{Row=25; New=36},
    @{Row=26; New=2650},
    @{Row=28; New=542},
    @{Row=30; New=329},
    @{Row=31; New=338},
    @{Row=32; New=389},
    @{Row=33; New=206},
    @{Row=34; New=24},
    @{Row=35; New=1598},
    @{Row=36; New=1002},
    @{Row=38; New=112},
    @{Row=39; New=72},
    @{Row=40; New=519},
    @{Row=41; New=493},
    @{Row=43; New=82},
    @{Row=44; New=610}
)

$wsExpo = $wb.Worksheets.Item("展览")

foreach ($c in $fChanges) {
    $wsExpo.Range("F" + $c.Row).Value = $c.New
    $wsAll.Range("F" + $c.Row).Value = $c.New
}
